# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.343.01'
$ws.Range("E2").Value = '  +2.40%  '
$ws.Range("D3").Value = '2.006.92'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.70'
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("E6").Value = '  +2.90%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.40'
$ws.Range("E7").Value = '  +4.23%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.383'
$ws.Range("E9").Value = '  +1.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0805'
$ws.Range("E10").Value = '  +2.17%  '
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.17'
$ws.Range("E12").Value = '  +10.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.73'
$ws.Range("E13").Value = '  +7.05%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.849'
$ws.Range("E14").Value = '  +2.83%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.299.97'
$ws.Range("E15").Value = '  +2.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.46'
$ws.Range("E16").Value = '  +3.29%  '
$ws.Range("D17").Value = '2.006.77'
$ws.Range("E17").Value = '  +3.10%  '
$ws.Range("D18").Value = '37.270.85'
$ws.Range("E18").Value = '  +2.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.51'
$ws.Range("E19").Value = '  +1.77%  '
$ws.Range("D20").Value = '0.0₃0868'
$ws.Range("E20").Value = '  +2.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.21'
$ws.Range("E21").Value = '  +3.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.15'
$ws.Range("E22").Value = '  +1.30%  '
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.50'
$ws.Range("E24").Value = '  +0.98%  '
$ws.Range("E25").Value = '  -0.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.146'
$ws.Range("E26").Value = '  +7.03%  '
$ws.Range("E27").Value = '  +3.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '163.92'
$ws.Range("E28").Value = '  +2.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.70'
$ws.Range("E29").Value = '  +2.18%  '
$ws.Range("E30").Value = '  +13.61%  '
$ws.Range("E31").Value = '  +1.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.83'
$ws.Range("E32").Value = '  +3.53%  '
$ws.Range("E33").Value = '  +7.11%  '
$ws.Range("E34").Value = '  +5.08%  '
$ws.Range("E35").Value = '  +6.21%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  +2.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.28'
$ws.Range("E38").Value = '  -3.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.53'
$ws.Range("E39").Value = '  +5.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0979'
$ws.Range("E40").Value = '  +0.79%  '
$ws.Range("E41").Value = '  +1.02%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.19'
$ws.Range("E42").Value = '  +2.69%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0215'
$ws.Range("E43").Value = '  +2.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.74'
$ws.Range("E44").Value = '  +6.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.77'
$ws.Range("E45").Value = '  +3.96%  '
$ws.Range("D46").Value = '1.374.15'
$ws.Range("E46").Value = '  +0.97%  '
$ws.Range("E47").Value = '  +2.83%  '
$ws.Range("E48").Value = '  +2.41%  '
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("E50").Value = '  +15.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.15'
$ws.Range("E51").Value = '  +5.62%  '
